$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-07-20 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-07-21 Sunday", 2) | Out-Null
$d.Content.Find.Execute("127÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "144÷6=", 2) | Out-Null
$d.Content.Find.Execute("576÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "150÷9=", 2) | Out-Null
$d.Content.Find.Execute("501÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "426÷4=", 2) | Out-Null
$d.Content.Find.Execute("764÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "776÷2=", 2) | Out-Null
$d.Content.Find.Execute("326÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "316÷7=", 2) | Out-Null
$d.Content.Find.Execute("638÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "924÷4=", 2) | Out-Null
$d.Content.Find.Execute("126÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "899÷7=", 2) | Out-Null
$d.Content.Find.Execute("892÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "482÷3=", 2) | Out-Null
$d.Content.Find.Execute("107÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "620÷4=", 2) | Out-Null
$d.Content.Find.Execute("759÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "880÷7=", 2) | Out-Null
$d.Content.Find.Execute("245÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "592÷8=", 2) | Out-Null
$d.Content.Find.Execute("127÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "756÷2=", 2) | Out-Null
$d.Content.Find.Execute("547÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "847÷7=", 2) | Out-Null
$d.Content.Find.Execute("819÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "368÷8=", 2) | Out-Null
$d.Content.Find.Execute("499÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "769÷8=", 2) | Out-Null
$d.Content.Find.Execute("702÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "224÷7=", 2) | Out-Null
$d.Content.Find.Execute("807÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "118÷3=", 2) | Out-Null
$d.Content.Find.Execute("509÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "851÷9=", 2) | Out-Null
$d.Content.Find.Execute("955÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "114÷8=", 2) | Out-Null
$d.Content.Find.Execute("815÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "637÷6=", 2) | Out-Null
$d.Content.Find.Execute("197÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "922÷4=", 2) | Out-Null
$d.Content.Find.Execute("307÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "203÷8=", 2) | Out-Null
$d.Content.Find.Execute("572÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "425÷5=", 2) | Out-Null
$d.Content.Find.Execute("646÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "182÷7=", 2) | Out-Null
$d.Content.Find.Execute("519÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "605÷9=", 2) | Out-Null
